# Generate Report for Archive
#
# Updates the localization status from "Ready for handoff" to "In
# Translation" on every sheet, and narrows the now-too-wide status/locale
# columns back down to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" -----------------

# Overview sheet: the zh-cn (E) and de-de (F) status columns, rows 2-3
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# zh-cn / de-de sheets: the Status column (C), rows 2-3
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Column widths: narrow the status/locale columns -----------------------
# (was 17.2159881591797 characters, now ~13.41 characters)

$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
